# Fruta / hortaliza, semanal
# Insert two new weekly price rows for "Perejil" at the top of the data
# block (rows 19-20), pushing the existing rows 19-30 down to 21-32.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 19; this shifts the old
# rows 19:30 down to 21:32 and extends the sheet dimension to A1:R32.
$ws.Rows("19:20").Insert()

# --- New row 19: Primera ---
$ws.Cells.Item(19, 1).Value = 7
$ws.Cells.Item(19, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(19, 3).Value = "Ñuble"
$ws.Cells.Item(19, 4).Value = 44883
$ws.Cells.Item(19, 5).Value = 16
$ws.Cells.Item(19, 6).Value = 100112044
$ws.Cells.Item(19, 7).Value = "Perejil"
$ws.Cells.Item(19, 8).Value = "Sin especificar"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 300
$ws.Cells.Item(19, 11).Value = 700
$ws.Cells.Item(19, 12).Value = 800
$ws.Cells.Item(19, 13).Value = 750
$ws.Cells.Item(19, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(19, 15).Value = "Región del Maule"
$ws.Cells.Item(19, 16).Value = 750
$ws.Cells.Item(19, 17).Value = 1
$ws.Cells.Item(19, 18).Value = "Hortaliza"

# --- New row 20: Segunda ---
$ws.Cells.Item(20, 1).Value = 7
$ws.Cells.Item(20, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(20, 3).Value = "Ñuble"
$ws.Cells.Item(20, 4).Value = 44883
$ws.Cells.Item(20, 5).Value = 16
$ws.Cells.Item(20, 6).Value = 100112044
$ws.Cells.Item(20, 7).Value = "Perejil"
$ws.Cells.Item(20, 8).Value = "Sin especificar"
$ws.Cells.Item(20, 9).Value = "Segunda"
$ws.Cells.Item(20, 10).Value = 200
$ws.Cells.Item(20, 11).Value = 600
$ws.Cells.Item(20, 12).Value = 600
$ws.Cells.Item(20, 13).Value = 600
$ws.Cells.Item(20, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(20, 15).Value = "Región del Maule"
$ws.Cells.Item(20, 16).Value = 600
$ws.Cells.Item(20, 17).Value = 1
$ws.Cells.Item(20, 18).Value = "Hortaliza"
